$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 10.26 = 42923.85 pesos`n✅ 42923.85 pesos = 10.19 = 946.02 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 97.44
$wsTasas.Range("O10").Value = 4182.5
$wsTasas.Range("N12").Value = 4212
$wsTasas.Range("O12").Value = 92.83
